$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.814.30'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '3.785.29'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.71%  '
$ws.Range("D7").Value = '3.781.07'
$ws.Range("E7").Value = '  -0.83%  '
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.78'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000246'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.97'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").Value = '4.419.38'
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("D16").Value = '3.797.47'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").Value = '67.827.13'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '458.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("E22").Value = '  -4.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.690'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("E24").Value = '  -1.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.79%  '
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("D30").Value = '3.925.88'
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("E31").Value = '  -6.55%  '
$ws.Range("E32").Value = '  -2.38%  '
$ws.Range("E33").Value = '  -1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.95'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.04%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("E38").Value = '  +7.21%  '
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.980'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.47%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.67%  '
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("E50").Value = '  -0.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.11%  '
